$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for all worker rows: 2507 -> 2508
$ws.Range("E16:E18").Value = "2508"
